$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.297.63"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "'2.507.82"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'321.69"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'108.25"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "'39.87"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").Value = "'20.31"
$ws.Range("E11").Value = "  +8.87%  "
$ws.Range("D12").Value = "'0.0819"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'7.19"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "'2.899.98"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "'2.506.23"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "'48.145.61"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "'13.13"
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("D20").Value = "'6.78"
$ws.Range("E20").Value = "  +2.10%  "
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("D22").Value = "'0.0₃0946"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = "'279.86"
$ws.Range("E23").Value = "  +13.20%  "
$ws.Range("D24").Value = "'72.29"
$ws.Range("E24").Value = "  +2.25%  "
$ws.Range("D25").Value = "'2.55"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'25.76"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "'2.21"
$ws.Range("E28").Value = "  -4.10%  "
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").Value = "'35.36"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("D33").Value = "'19.66"
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("D34").Value = "'5.35"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'1.01"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'0.0783"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").Value = "'4.66"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "'121.71"
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'21.54"
$ws.Range("E43").Value = "  -4.05%  "
$ws.Range("D44").Value = "'0.0303"
$ws.Range("E44").Value = "  +1.83%  "
$ws.Range("D45").Value = "'2.017.43"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").Value = "'3.17"
$ws.Range("E46").Value = "  +4.26%  "
$ws.Range("D49").Value = "'9.03"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Value = "'5.18"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("D51").Value = "'80.60"
$ws.Range("E51").Value = "  +3.70%  "

# Row 47/48: Stacks and ApeXProtocol swap positions with new volume values
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'1.99"
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.85"
$ws.Range("E48").Value = "  +3.00%  "
